$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.087.64"
$ws.Range("E2").Value = "  -1.35%  "

$ws.Range("D3").Value = "2.306.84"
$ws.Range("E3").Value = "  -2.15%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.21"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.10"
$ws.Range("E6").Value = "  -1.73%  "

$ws.Range("E7").Value = "  -1.34%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  -1.61%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.89"
$ws.Range("E10").Value = "  -3.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.49"
$ws.Range("E12").Value = "  +0.27%  "

$ws.Range("E13").Value = "  +0.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.981"
$ws.Range("E14").Value = "  -2.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.52"
$ws.Range("E15").Value = "  -3.03%  "

$ws.Range("D16").Value = "2.655.24"
$ws.Range("E16").Value = "  -2.10%  "

$ws.Range("D17").Value = "2.305.35"
$ws.Range("E17").Value = "  -4.09%  "

$ws.Range("D18").Value = "42.188.87"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.77"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.90"
$ws.Range("E21").Value = "  -3.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "281.67"
$ws.Range("E22").Value = "  +8.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.58"
$ws.Range("E23").Value = "  -0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.29"
$ws.Range("E24").Value = "  -0.82%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.98"
$ws.Range("E25").Value = "  +5.74%  "

$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.98"
$ws.Range("E27").Value = "  +0.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.99"
$ws.Range("E28").Value = "  -3.74%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "23.28"
$ws.Range("E29").Value = "  +0.96%  "

$ws.Range("E30").Value = "  +0.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "165.75"
$ws.Range("E31").Value = "  -5.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.53"
$ws.Range("E32").Value = "  -2.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0884"
$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("E34").Value = "  -1.11%  "

$ws.Range("E35").Value = "  -3.82%  "

$ws.Range("E36").Value = "  +1.66%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.119"
$ws.Range("E37").Value = "  -4.80%  "

$ws.Range("E38").Value = "  +1.23%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0353"
$ws.Range("E39").Value = "  -2.38%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.91"
$ws.Range("E40").Value = "  +7.98%  "

$ws.Range("E41").Value = "  -3.43%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.62"
$ws.Range("E42").Value = "  +19.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.50"
$ws.Range("E43").Value = "  +1.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.90"
$ws.Range("E44").Value = "  -1.13%  "

$ws.Range("E45").Value = "  -4.34%  "

$ws.Range("E46").Value = "  +0.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "117.35"

$ws.Range("E48").Value = "  +0.70%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "78.80"
$ws.Range("E49").Value = "  +6.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.17"
$ws.Range("E50").Value = "  +0.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.35"
$ws.Range("E51").Value = "  -2.60%  "
